$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.583803653717041
$ws.Range("B1").Value = 0.7355457544326782
$ws.Range("C1").Value = 1.064057350158691
$ws.Range("D1").Value = 4.302477836608887
$ws.Range("E1").Value = 3.983104228973389
